# Auto-generated edit script: updates market-price columns (H:N) on specific
# leve rows across all 8 worksheets, per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 99.2
$ws.Range("I9").Value = 99
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 99
$ws.Range("L9").Value = 100
$ws.Range("M9").Value = 70
$ws.Range("N9").Value = -438

# Row 106
$ws.Range("H106").Value = 71714970
$ws.Range("I106").Value = 400792.4
$ws.Range("J106").Value = 250000400
$ws.Range("K106").Value = 400792.4
$ws.Range("L106").Value = 250000400
$ws.Range("M106").Value = -400161.4
$ws.Range("N106").Value = -250001662

# Row 108
$ws.Range("H108").Value = 29163.637
$ws.Range("J108").Value = 29163.637
$ws.Range("L108").Value = 29163.637
$ws.Range("N108").Value = -36843.637

# Row 113
$ws.Range("H113").Value = 1961.9474
$ws.Range("I113").Value = 1587.1818
$ws.Range("J113").Value = 2477.25
$ws.Range("K113").Value = 1587.1818
$ws.Range("L113").Value = 2477.25
$ws.Range("M113").Value = 1666.8182
$ws.Range("N113").Value = -8985.25

# Row 129
$ws.Range("H129").Value = 13414.95
$ws.Range("J129").Value = 15473.275
$ws.Range("L129").Value = 46419.825
$ws.Range("N129").Value = -56419.825

$ws = $wb.Worksheets.Item("ARM")
# Row 92
$ws.Range("H92").Value = 37641.668
$ws.Range("J92").Value = 37641.668
$ws.Range("L92").Value = 37641.668
$ws.Range("N92").Value = -42633.668

# Row 122
$ws.Range("H122").Value = 3571
$ws.Range("I122").Value = 2012
$ws.Range("J122").Value = 3882.8
$ws.Range("K122").Value = 6036
$ws.Range("L122").Value = 11648.4
$ws.Range("M122").Value = -3586
$ws.Range("N122").Value = -16548.4

# Row 132
$ws.Range("H132").Value = 3938.8809
$ws.Range("I132").Value = 5356.5415
$ws.Range("J132").Value = 2048.6667
$ws.Range("K132").Value = 16069.6245
$ws.Range("L132").Value = 6146.000100000001
$ws.Range("M132").Value = -13539.6245
$ws.Range("N132").Value = -11206.0001

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 2804.5
$ws.Range("I107").Value = 2804.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2804.5
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -884.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 9092689
$ws.Range("I31").Value = 1914.7059
$ws.Range("J31").Value = 40001320
$ws.Range("K31").Value = 1914.7059
$ws.Range("L31").Value = 40001320
$ws.Range("M31").Value = -1619.7059
$ws.Range("N31").Value = -40001910

# Row 34
$ws.Range("H34").Value = 9092689
$ws.Range("I34").Value = 1914.7059
$ws.Range("J34").Value = 40001320
$ws.Range("K34").Value = 1914.7059
$ws.Range("L34").Value = 40001320
$ws.Range("M34").Value = -1712.7059
$ws.Range("N34").Value = -40001724

# Row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").ClearContents()
$ws.Range("N108").Value = 0

$ws = $wb.Worksheets.Item("CUL")
# Row 70
$ws.Range("H70").Value = 5253
$ws.Range("I70").Value = 1012
$ws.Range("J70").Value = 6666.6665
$ws.Range("K70").Value = 3036
$ws.Range("L70").Value = 19999.9995
$ws.Range("M70").Value = -2721
$ws.Range("N70").Value = -20629.9995

# Row 73
$ws.Range("H73").Value = 5253
$ws.Range("I73").Value = 1012
$ws.Range("J73").Value = 6666.6665
$ws.Range("K73").Value = 3036
$ws.Range("L73").Value = 19999.9995
$ws.Range("M73").Value = -1944
$ws.Range("N73").Value = -22183.9995

# Row 117
$ws.Range("H117").Value = 3997.1428
$ws.Range("J117").Value = 3997.1428
$ws.Range("L117").Value = 11991.4284
$ws.Range("N117").Value = -18875.4284

# Row 131
$ws.Range("H131").Value = 16323.286
$ws.Range("I131").Value = 18905
$ws.Range("J131").Value = 833
$ws.Range("K131").Value = 56715
$ws.Range("L131").Value = 2499
$ws.Range("M131").Value = -51675
$ws.Range("N131").Value = -12579

# Row 137
$ws.Range("H137").Value = 45482160
$ws.Range("I137").Value = 6110
$ws.Range("J137").Value = 55587950
$ws.Range("K137").Value = 18330
$ws.Range("L137").Value = 166763850
$ws.Range("M137").Value = -13230
$ws.Range("N137").Value = -166774050

$ws = $wb.Worksheets.Item("GSM")
# Row 20
$ws.Range("H20").Value = 48000
$ws.Range("J20").Value = 48000
$ws.Range("L20").Value = 48000
$ws.Range("N20").Value = -48490

# Row 24
$ws.Range("H24").Value = 30000000
$ws.Range("I24").Value = 30000000
$ws.Range("K24").Value = 30000000
$ws.Range("M24").Value = -29999827

# Row 113
$ws.Range("H113").Value = 62502100
$ws.Range("I113").Value = 83335660
$ws.Range("J113").Value = 1400
$ws.Range("K113").Value = 83335660
$ws.Range("L113").Value = 1400
$ws.Range("M113").Value = -83333490
$ws.Range("N113").Value = -5740

# Row 122
$ws.Range("H122").Value = 2615.2856
$ws.Range("I122").Value = 2185.6667
$ws.Range("K122").Value = 6557.000100000001
$ws.Range("M122").Value = -4107.000100000001

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 1500
$ws.Range("I61").Value = 1500
$ws.Range("K61").Value = 1500
$ws.Range("M61").Value = -1298

# Row 113
$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 1500
$ws.Range("K113").Value = 1500
$ws.Range("M113").Value = 670

# Row 122
$ws.Range("H122").Value = 2554.353
$ws.Range("I122").Value = 2432.6155
$ws.Range("J122").Value = 2950
$ws.Range("K122").Value = 7297.8465
$ws.Range("L122").Value = 8850
$ws.Range("M122").Value = -4847.8465
$ws.Range("N122").Value = -13750

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 290.7143
$ws.Range("I107").Value = 290.7143
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 872.1428999999999
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = 1047.8571

# Row 114
$ws.Range("H114").Value = 29875
$ws.Range("J114").Value = 29875
$ws.Range("L114").Value = 29875
$ws.Range("N114").Value = -38553

# Row 122
$ws.Range("H122").Value = 1327.7273
$ws.Range("I122").Value = 1377.6875
$ws.Range("J122").Value = 1194.5
$ws.Range("K122").Value = 4133.0625
$ws.Range("L122").Value = 3583.5
$ws.Range("M122").Value = -1683.0625
$ws.Range("N122").Value = -8483.5

# Row 138
$ws.Range("H138").Value = 111214.5
$ws.Range("J138").Value = 111214.5
$ws.Range("L138").Value = 111214.5
$ws.Range("N138").Value = -121494.5
